# Updated cryptos list on Sat Jun 29 17:40:07 UTC 2024 with GitHub Actions
#
# The Price column (D) holds plain numeric-looking text (Excel would
# auto-convert "571.48" etc. to a real number on assignment), so those
# writes are apostrophe-prefixed to force text, matching the original
# inlineStr/text storage. The Volume column (E) always carries padding
# spaces + a percent sign, so it never round-trips as a number anyway.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'60.957.25"
$ws.Range("E2").Value = "  -0.01%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'3.384.77"
$ws.Range("E3").Value = "  -0.25%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.02%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'571.48"
$ws.Range("E5").Value = "  -0.19%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'142.01"
$ws.Range("E6").Value = "  -0.45%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.06%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +0.25%  "

# Row 9 - Toncoin
$ws.Range("D9").Value = "'7.63"
$ws.Range("E9").Value = "  +1.37%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -1.26%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  -1.36%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "'3.963.63"
$ws.Range("E12").Value = "  -0.27%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +1.91%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "'27.86"
$ws.Range("E14").Value = "  -1.43%  "

# Row 15 - was WrappedEther, now ShibaInu (rows 15/16 swapped)
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "'0.0000172"
$ws.Range("E15").Value = "  +0.15%  "

# Row 16 - was ShibaInu, now WrappedEther
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "'3.399.73"
$ws.Range("E16").Value = "  -0.07%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "'61.075.29"
$ws.Range("E17").Value = "  +0.03%  "

# Row 18 - Polkadot
$ws.Range("D18").Value = "'6.09"
$ws.Range("E18").Value = "  -3.21%  "

# Row 19 - Chainlink
$ws.Range("E19").Value = "  -4.00%  "

# Row 20 - Uniswap
$ws.Range("D20").Value = "'8.94"
$ws.Range("E20").Value = "  -1.99%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "'384.75"
$ws.Range("E21").Value = "  -1.16%  "

# Row 22 - Litecoin
$ws.Range("D22").Value = "'75.06"
$ws.Range("E22").Value = "  +2.48%  "

# Row 23 - Polygon
$ws.Range("D23").Value = "'0.553"
$ws.Range("E23").Value = "  -1.68%  "

# Row 25 - PEPE
$ws.Range("E25").Value = "  -3.28%  "

# Row 26 - WrappedeETH
$ws.Range("D26").Value = "'3.522.17"
$ws.Range("E26").Value = "  -0.36%  "

# Row 27 - Kaspa
$ws.Range("E27").Value = "  +1.44%  "

# Row 28 - Binance-PegBSC-USD
$ws.Range("E28").Value = "  -0.04%  "

# Row 29 - RenderToken
$ws.Range("E29").Value = "  -1.96%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  -0.39%  "

# Row 31 - InternetComputer(DFINITY)
$ws.Range("D31").Value = "'7.96"
$ws.Range("E31").Value = "  -2.63%  "

# Row 32 - USDe
$ws.Range("E32").Value = "  -0.03%  "

# Row 33 - Fetch.AI
$ws.Range("E33").Value = "  -4.46%  "

# Row 34 - EthereumClassic
$ws.Range("E34").Value = "  -2.72%  "

# Row 35 - Aptos
$ws.Range("D35").Value = "'6.95"
$ws.Range("E35").Value = "  -0.25%  "

# Row 36 - Monero
$ws.Range("D36").Value = "'166.94"
$ws.Range("E36").Value = "  -0.05%  "

# Row 37 - RenzoRestakedETH
$ws.Range("D37").Value = "'3.417.37"
$ws.Range("E37").Value = "  -0.09%  "

# Row 38 - NEARProtocol
$ws.Range("E38").Value = "  -2.02%  "

# Row 39 - ImmutableX
$ws.Range("E39").Value = "  -3.81%  "

# Row 40 - Hedera
$ws.Range("D40").Value = "'0.0768"
$ws.Range("E40").Value = "  -1.90%  "

# Row 41 - EnergySwap
$ws.Range("D41").Value = "'26.91"
$ws.Range("E41").Value = "  -0.07%  "

# Row 42 - FirstDigitalUSD
$ws.Range("E42").Value = "  +0.00%  "

# Row 43 - Mantle
$ws.Range("E43").Value = "  -0.99%  "

# Row 44 - Filecoin
$ws.Range("D44").Value = "'4.37"
$ws.Range("E44").Value = "  -2.15%  "

# Row 45 - Stacks
$ws.Range("E45").Value = "  -2.36%  "

# Row 46 - ONDO
$ws.Range("E46").Value = "  -0.82%  "

# Row 47 - Maker
$ws.Range("D47").Value = "'2.451.86"
$ws.Range("E47").Value = "  -3.57%  "

# Row 48 - InjectiveProtocol
$ws.Range("D48").Value = "'22.92"
$ws.Range("E48").Value = "  -0.12%  "

# Row 49 - Cosmos
$ws.Range("D49").Value = "'6.71"
$ws.Range("E49").Value = "  -2.42%  "

# Row 50 - dogwifhat
$ws.Range("D50").Value = "'2.16"
$ws.Range("E50").Value = "  +9.83%  "

# Row 51 - VeChain
$ws.Range("E51").Value = "  +1.59%  "
